$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 with the new feature/tags/scenario values
$ws.Range("A2").Value = "D:\IBK\Proyectos\documnetadorfeatures\test2\Login.feature"
$ws.Range("B2").Value = "['@Automated', '@Happy_path', '@Functional_testing', '@MilesRegression', '@regressionTest', '@LoginTest']"
$ws.Range("C2").Value = "[HAPPY PATH] Validar el ingreso con diferentes tipo de documento - Usuario existente"

# Add a new row 3 for the second scenario of the same feature
$ws.Range("A3").Value = "D:\IBK\Proyectos\documnetadorfeatures\test2\Login.feature"
$ws.Range("B3").Value = "['@Automated', '@Happy_path', '@Functional_testing', '@MilesRegression']"
$ws.Range("C3").Value = "[HAPPY PATH] Validar el ingreso a Mi cuenta"
